$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helpers: copy a cell's value and/or formatting without creating new
# style entries for attributes that are already at their default value
# (Excel's xlsx writer only allocates a new <xf> when a property differs
# from the implicit default, so touching a default-valued property can
# otherwise spawn a duplicate style).
# ---------------------------------------------------------------------
function Copy-CellStyle {
    param($srcRange, $dstRange)
    if ($srcRange.NumberFormat -ne "General") {
        $dstRange.NumberFormat = $srcRange.NumberFormat
    }
    if ($srcRange.Font.Bold) {
        $dstRange.Font.Bold = $true
    }
    if ($srcRange.HorizontalAlignment -ne 1) {
        $dstRange.HorizontalAlignment = $srcRange.HorizontalAlignment
    }
}

function Copy-Cell {
    param($srcRange, $dstRange)
    $dstRange.Value = $srcRange.Value2
    Copy-CellStyle $srcRange $dstRange
}

# ---------------------------------------------------------------------
# Sheet "Key Metrics": add a Q4FY22 column (D) mirroring column C.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Key Metrics")

$ws1.Range("D1").Value = "Q4FY22"
Copy-CellStyle $ws1.Range("C1") $ws1.Range("D1")

2..12 | ForEach-Object {
    $row = $_
    Copy-Cell $ws1.Range("C$row") $ws1.Range("D$row")
}

$ws1.Range("B1:D1").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Balance Sheet": add Q4FY24 / Q4FY23 / Q4FY22 columns (E/F/G)
# mirroring columns B/C/C for both tables (rows 1-8 and rows 11-18).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Balance Sheet")

# Header rows get the same quarter labels (and bold/centered style) as
# the "Key Metrics" sheet's header row -- not a copy of this sheet's own
# B/C headers (which read "FY23"/"FY24", a different label set).
$headerRows = 1, 11
foreach ($row in $headerRows) {
    $ws3.Range("E$row").Value = "Q4FY24"
    $ws3.Range("F$row").Value = "Q4FY23"
    $ws3.Range("G$row").Value = "Q4FY22"
    Copy-CellStyle $ws1.Range("B1") $ws3.Range("E$row")
    Copy-CellStyle $ws1.Range("C1") $ws3.Range("F$row")
    Copy-CellStyle $ws1.Range("D1") $ws3.Range("G$row")
}

$dataRows = 2, 3, 4, 5, 6, 7, 8, 12, 13, 14, 15, 16, 17, 18
foreach ($row in $dataRows) {
    Copy-Cell $ws3.Range("B$row") $ws3.Range("E$row")
    Copy-Cell $ws3.Range("C$row") $ws3.Range("F$row")
    Copy-Cell $ws3.Range("C$row") $ws3.Range("G$row")
}

$ws3.Range("E11:G11").Select() | Out-Null

# Re-activate the Balance Sheet tab (it was active/selected before the
# edit, and selecting ranges on other sheets above would otherwise have
# switched the active tab away from it).
$ws3.Activate() | Out-Null
